$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.712061882019043
$ws.Range("B1").Value = 2.859842538833618
$ws.Range("C1").Value = 3.557749271392822
$ws.Range("D1").Value = 1.378337264060974
$ws.Range("E1").Value = 0.921788215637207
